# Apply edits from covid_disparities_output_2020-07-08.xlsx commit
# "Results from July 08, 2020 06:13:49 PM America/Chicago TZ run"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = $False
$ws.Range("O2").Value = "An error occurred. ... KeyError(`"None of ['Race/Ethnicity'] are in the columns`")"

# Row 4
$ws.Range("B4").Value = 44020
$ws.Range("C4").Value = "'214570"
$ws.Range("D4").Value = "'18618"
$ws.Range("E4").Value = 33320
$ws.Range("F4").Value = 5210
$ws.Range("H4").Value = 30.52
$ws.Range("K4").Value = 110482
$ws.Range("L4").Value = 17070

# Row 6
$ws.Range("B6").Value = 44020
$ws.Range("C6").Value = 55986
$ws.Range("D6").Value = 685
$ws.Range("E6").Value = 11560
$ws.Range("F6").Value = 241
$ws.Range("G6").Value = 20.65
$ws.Range("H6").Value = 35.18

# Row 7
$ws.Range("B7").Value = 44020
$ws.Range("C7").Value = "'26755"
$ws.Range("D7").Value = "'201"
$ws.Range("E7").Value = "'689"

# Row 8
$ws.Range("B8").Value = 44020
$ws.Range("C8").Value = 17919
$ws.Range("D8").Value = 608
$ws.Range("E8").Value = 1707
$ws.Range("G8").Value = 13.77
$ws.Range("H8").Value = 15.41
$ws.Range("K8").Value = 12399
$ws.Range("L8").Value = 570

# Row 9
$ws.Range("B9").Value = 44020

# Row 12
$ws.Range("B12").Value = 44020
$ws.Range("C12").Value = 14017
$ws.Range("D12").Value = 527
$ws.Range("E12").Value = 266

# Row 16
$ws.Range("B16").Value = 44019
$ws.Range("C16").Value = 123004
$ws.Range("D16").Value = 3642
$ws.Range("E16").Value = 3310
$ws.Range("F16").Value = 372
$ws.Range("G16").Value = 4.74
$ws.Range("H16").Value = 10.98
$ws.Range("K16").Value = 69866
$ws.Range("L16").Value = 3389

# Row 18
$ws.Range("B18").Value = 44019
$ws.Range("C18").Value = 32888
$ws.Range("D18").Value = 1188
$ws.Range("E18").Value = 15720
$ws.Range("F18").Value = 595
$ws.Range("G18").Value = 47.8
$ws.Range("H18").Value = 50.08

# Row 23
$ws.Range("B23").Value = 44020
$ws.Range("C23").Value = 35116
$ws.Range("D23").Value = 1704
$ws.Range("E23").Value = 1840
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 6.39
$ws.Range("H23").Value = 6.75
$ws.Range("K23").Value = 28809
$ws.Range("L23").Value = 1644

# Row 24
$ws.Range("B24").Value = 44020
$ws.Range("C24").Value = 20425
$ws.Range("E24").Value = 1205
$ws.Range("G24").Value = 7.63
$ws.Range("K24").Value = 15786

# Row 25
$ws.Range("B25").Value = 44020
$ws.Range("C25").Value = 67153
$ws.Range("D25").Value = 5934
$ws.Range("E25").Value = 20210
$ws.Range("F25").Value = 2367
$ws.Range("G25").Value = 30.1

# Row 28
$ws.Range("B28").Value = 44020
$ws.Range("C28").Value = 1226
$ws.Range("E28").Value = 30
$ws.Range("G28").Value = 1.41
$ws.Range("K28").Value = 2135

# Row 29
$ws.Range("B29").Value = 44020
$ws.Range("C29").Value = 33154
$ws.Range("D29").Value = 807
$ws.Range("E29").Value = 5634
$ws.Range("G29").Value = 18.88
$ws.Range("H29").Value = 24.15
$ws.Range("K29").Value = 29848
$ws.Range("L29").Value = 795

# Row 30
$ws.Range("B30").Value = 44020
$ws.Range("C30").Value = 103890
$ws.Range("D30").Value = 2922
$ws.Range("E30").Value = 28382
$ws.Range("F30").Value = 1372
$ws.Range("G30").Value = 27.32

# Row 31
$ws.Range("B31").Value = 44020
$ws.Range("C31").Value = 37941
$ws.Range("D31").Value = 1394
$ws.Range("E31").Value = 1478
$ws.Range("H31").Value = 3.4
$ws.Range("K31").Value = 26965
$ws.Range("L31").Value = 1323

# Row 33
$ws.Range("B33").Value = 44020
$ws.Range("C33").Value = 12462
$ws.Range("D33").Value = 515
$ws.Range("E33").Value = 3197
$ws.Range("G33").Value = 25.65
$ws.Range("H33").Value = 25.44

# Row 37
$ws.Range("B37").Value = 44020
$ws.Range("C37").Value = 149432
$ws.Range("D37").Value = 7099
$ws.Range("E37").Value = 25072
$ws.Range("F37").Value = 1967
$ws.Range("G37").Value = 16.78
$ws.Range("H37").Value = 27.71

# Row 38
$ws.Range("B38").Value = 44020
$ws.Range("C38").Value = 8969
$ws.Range("D38").Value = 98
$ws.Range("E38").Value = 134
$ws.Range("H38").Value = 1.02

# Row 40
$ws.Range("B40").Value = 44020
$ws.Range("C40").Value = 110602
$ws.Range("D40").Value = 8243
$ws.Range("E40").Value = 10414
$ws.Range("F40").Value = 675
$ws.Range("G40").Value = 9.42
$ws.Range("H40").Value = 8.19

# Row 41
$ws.Range("B41").Value = 44020
$ws.Range("C41").Value = 12834
$ws.Range("E41").Value = 3731
$ws.Range("G41").Value = 31.32
$ws.Range("K41").Value = 11914
